$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "product_family" column (E) and the "upstream_processes" /
# "downstream_processes" columns (which, after deleting E, land in H:I).
# Delete from right to left so earlier deletions don't shift the
# still-to-be-processed column references.
$ws.Range("I1:J1").EntireColumn.Delete()
$ws.Range("E1").EntireColumn.Delete()
